# Daily price-data refresh: insert a new "today" row at the top of the
# table (row 2, just below the header) with the next day's date and the
# same commodity prices, pushing all the historical rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new day by inserting a row above the current top
# data row; everything below (including the old row 2) shifts down one.
$ws.Range("A2").EntireRow.Insert()

# Force column A to be stored as literal text (matching every other date
# cell in the sheet) instead of letting Excel auto-convert the
# date-shaped string into a date serial number. Clear the formatting
# afterwards so the cell's style index stays the default (same as its
# neighbours) rather than keeping a lingering "@" text format.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2026-02-19"
$ws.Range("A2").ClearFormats()

# Carry forward the same (unchanged) commodity prices used throughout
# the rest of the sheet.
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
